$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# OS sheet: enter the missing "final" marks obtained for G9 (was blank)
# ---------------------------------------------------------------------
$os = $wb.Worksheets.Item("OS")
$os.Range("G9").Value = 210

# ---------------------------------------------------------------------
# POM sheet: update Project / Quiz / Assignment marks
# ---------------------------------------------------------------------
$pom = $wb.Worksheets.Item("POM")

# Project row (7)
$pom.Range("C7").Value = 20
$pom.Range("D7").Value = 20
$pom.Range("G7").Value = 20

# Quiz row (8)
$pom.Range("G8").Value = 8
$pom.Range("H8").Value = 7
$pom.Range("J8").Value = 7
$pom.Range("K8").Value = 6

# Assignment row (9)
$pom.Range("D9").Value = 50
$pom.Range("J9").Value = 10
# K9 is a brand-new cell in this row; copy the formatting from J9 first
# so it keeps the same direct style (s="1") as its neighbours.
$pom.Range("J9").Copy()
$pom.Range("K9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$pom.Range("K9").Value = 10

# CP row (12) is removed entirely
$pom.Range("B12:G12").ClearContents()

# Restore the selection on POM to match the authored state
$pom.Range("F11").Select()

# ---------------------------------------------------------------------
# OS sheet selection
# ---------------------------------------------------------------------
$os.Range("G10").Select()

# ---------------------------------------------------------------------
# GPA sheet: update the selection/scroll state
# ---------------------------------------------------------------------
$gpa = $wb.Worksheets.Item("GPA")
$gpa.Activate()
$gpa.Range("D10:I10").Select()
